$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update forecasted values in column F (MSE) per new model run
$ws.Range("F2").Value = 7.7214
$ws.Range("F3").Value = 7.7594000000000003
$ws.Range("F4").Value = 7.7618
$ws.Range("F6").Value = 5.9039999999999999
$ws.Range("F7").Value = 5.9157999999999999
$ws.Range("F8").Value = 5.9353999999999996

# Move selection to F9 as last active cell
$ws.Activate()
$ws.Range("F9").Select()
